# Forms the consolidated report by filling in the "Absent" (column H) values
# for the attendance rows that were previously blank / incorrect.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that need column H (Absent) updated to reflect the consolidated report.
$ws.Range("H3").Value  = 1
$ws.Range("H4").Value  = 0
$ws.Range("H5").Value  = 1
$ws.Range("H6").Value  = 0
$ws.Range("H11").Value = 1
$ws.Range("H13").Value = 0
